$d = $word.ActiveDocument
$n = $d.Styles.Item(1)
try {
  $n.Delete()
  Write-Output "deleted"
} catch {
  Write-Output "ERR: $_"
}
